# Auto-generated Excel COM-interop edit script
# Applies the Twintania_Profits market-data refresh described in the commit diff.
# For every changed cell we assert the previous (pre-edit) value matches what the
# diff recorded, then write the new value -- this guards against mis-mapped rows.

$wb = $excel.ActiveWorkbook
$mismatchCount = 0

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 13
$cell = $ws.Cells.Item(13, 8)
if ($cell.Value2 -ne 16070) { $mismatchCount++ }
$cell.Value = 9600
$cell = $ws.Cells.Item(13, 10)
if ($cell.Value2 -ne 17581.666) { $mismatchCount++ }
$cell.Value = 10900
$cell = $ws.Cells.Item(13, 12)
if ($cell.Value2 -ne 17581.666) { $mismatchCount++ }
$cell.Value = 10900
$cell = $ws.Cells.Item(13, 14)
if ($cell.Value2 -ne -17919.666) { $mismatchCount++ }
$cell.Value = -11238

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 10
$cell = $ws.Cells.Item(10, 8)
if ($cell.Value2 -ne 3622.1) { $mismatchCount++ }
$cell.Value = 4279
$cell = $ws.Cells.Item(10, 10)
if ($cell.Value2 -ne 889.6667) { $mismatchCount++ }
$cell.Value = 837.25
$cell = $ws.Cells.Item(10, 12)
if ($cell.Value2 -ne 889.6667) { $mismatchCount++ }
$cell.Value = 837.25
$cell = $ws.Cells.Item(10, 14)
if ($cell.Value2 -ne -1229.6667) { $mismatchCount++ }
$cell.Value = -1177.25
# Row 12
$cell = $ws.Cells.Item(12, 8)
if ($cell.Value2 -ne 2883.818) { $mismatchCount++ }
$cell.Value = 2450.7778
$cell = $ws.Cells.Item(12, 9)
if ($cell.Value2 -ne 4281) { $mismatchCount++ }
$cell.Value = 3486
$cell = $ws.Cells.Item(12, 10)
if ($cell.Value2 -ne 2085.4285) { $mismatchCount++ }
$cell.Value = 1933.1666
$cell = $ws.Cells.Item(12, 11)
if ($cell.Value2 -ne 4281) { $mismatchCount++ }
$cell.Value = 3486
$cell = $ws.Cells.Item(12, 12)
if ($cell.Value2 -ne 2085.4285) { $mismatchCount++ }
$cell.Value = 1933.1666
$cell = $ws.Cells.Item(12, 13)
if ($cell.Value2 -ne -4108) { $mismatchCount++ }
$cell.Value = -3313
$cell = $ws.Cells.Item(12, 14)
if ($cell.Value2 -ne -2431.4285) { $mismatchCount++ }
$cell.Value = -2279.1666
# Row 29
$cell = $ws.Cells.Item(29, 8)
if ($cell.Value2 -ne 19989.5) { $mismatchCount++ }
$cell.Value = 10994.5
$cell = $ws.Cells.Item(29, 10)
if ($cell.Value2 -ne 19989.5) { $mismatchCount++ }
$cell.Value = 10994.5
$cell = $ws.Cells.Item(29, 12)
if ($cell.Value2 -ne 19989.5) { $mismatchCount++ }
$cell.Value = 10994.5
$cell = $ws.Cells.Item(29, 14)
if ($cell.Value2 -ne -20605.5) { $mismatchCount++ }
$cell.Value = -11610.5
# Row 32
$cell = $ws.Cells.Item(32, 8)
if ($cell.Value2 -ne 3571.1736) { $mismatchCount++ }
$cell.Value = 3539.192
$cell = $ws.Cells.Item(32, 9)
if ($cell.Value2 -ne 1766.5333) { $mismatchCount++ }
$cell.Value = 1751.5714
$cell = $ws.Cells.Item(32, 11)
if ($cell.Value2 -ne 1766.5333) { $mismatchCount++ }
$cell.Value = 1751.5714
$cell = $ws.Cells.Item(32, 13)
if ($cell.Value2 -ne -1479.5333) { $mismatchCount++ }
$cell.Value = -1464.5714
# Row 50
$cell = $ws.Cells.Item(50, 8)
if ($cell.Value2 -ne 891) { $mismatchCount++ }
$cell.Value = 907.6667
$cell = $ws.Cells.Item(50, 9)
if ($cell.Value2 -ne 211.25) { $mismatchCount++ }
$cell.Value = 236
$cell = $ws.Cells.Item(50, 10)
if ($cell.Value2 -ne 2250.5) { $mismatchCount++ }
$cell.Value = 2251
$cell = $ws.Cells.Item(50, 11)
if ($cell.Value2 -ne 211.25) { $mismatchCount++ }
$cell.Value = 236
$cell = $ws.Cells.Item(50, 12)
if ($cell.Value2 -ne 2250.5) { $mismatchCount++ }
$cell.Value = 2251
$cell = $ws.Cells.Item(50, 13)
if ($cell.Value2 -ne 502.75) { $mismatchCount++ }
$cell.Value = 478
$cell = $ws.Cells.Item(50, 14)
if ($cell.Value2 -ne -3678.5) { $mismatchCount++ }
$cell.Value = -3679

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$cell = $ws.Cells.Item(5, 8)
if ($cell.Value2 -ne 1625.4) { $mismatchCount++ }
$cell.Value = 1549.6
$cell = $ws.Cells.Item(5, 9)
if ($cell.Value2 -ne 664) { $mismatchCount++ }
$cell.Value = 649.3333
$cell = $ws.Cells.Item(5, 10)
if ($cell.Value2 -ne 2266.3333) { $mismatchCount++ }
$cell.Value = 2900
$cell = $ws.Cells.Item(5, 11)
if ($cell.Value2 -ne 664) { $mismatchCount++ }
$cell.Value = 649.3333
$cell = $ws.Cells.Item(5, 12)
if ($cell.Value2 -ne 2266.3333) { $mismatchCount++ }
$cell.Value = 2900
$cell = $ws.Cells.Item(5, 13)
if ($cell.Value2 -ne -551) { $mismatchCount++ }
$cell.Value = -536.3333
$cell = $ws.Cells.Item(5, 14)
if ($cell.Value2 -ne -2492.3333) { $mismatchCount++ }
$cell.Value = -3126
# Row 11
$cell = $ws.Cells.Item(11, 8)
if ($cell.Value2 -ne 1569.4286) { $mismatchCount++ }
$cell.Value = 362
$cell = $ws.Cells.Item(11, 9)
if ($cell.Value2 -ne 400) { $mismatchCount++ }
$cell.Value = 71.46154
$cell = $ws.Cells.Item(11, 10)
if ($cell.Value2 -ne 1764.3334) { $mismatchCount++ }
$cell.Value = 1117.4
$cell = $ws.Cells.Item(11, 11)
if ($cell.Value2 -ne 400) { $mismatchCount++ }
$cell.Value = 71.46154
$cell = $ws.Cells.Item(11, 12)
if ($cell.Value2 -ne 1764.3334) { $mismatchCount++ }
$cell.Value = 1117.4
$cell = $ws.Cells.Item(11, 13)
if ($cell.Value2 -ne -260) { $mismatchCount++ }
$cell.Value = 68.53846
$cell = $ws.Cells.Item(11, 14)
if ($cell.Value2 -ne -2044.3334) { $mismatchCount++ }
$cell.Value = -1397.4
# Row 12
$cell = $ws.Cells.Item(12, 8)
if ($cell.Value2 -ne 2199.2856) { $mismatchCount++ }
$cell.Value = 1098.4286
$cell = $ws.Cells.Item(12, 9)
if ($cell.Value2 -ne 2599) { $mismatchCount++ }
$cell.Value = 1448.25
$cell = $ws.Cells.Item(12, 10)
if ($cell.Value2 -ne 1666.3334) { $mismatchCount++ }
$cell.Value = 632
$cell = $ws.Cells.Item(12, 11)
if ($cell.Value2 -ne 2599) { $mismatchCount++ }
$cell.Value = 1448.25
$cell = $ws.Cells.Item(12, 12)
if ($cell.Value2 -ne 1666.3334) { $mismatchCount++ }
$cell.Value = 632
$cell = $ws.Cells.Item(12, 13)
if ($cell.Value2 -ne -2431) { $mismatchCount++ }
$cell.Value = -1280.25
$cell = $ws.Cells.Item(12, 14)
if ($cell.Value2 -ne -2002.3334) { $mismatchCount++ }
$cell.Value = -968
# Row 23
$cell = $ws.Cells.Item(23, 8)
if ($cell.Value2 -ne 3515.25) { $mismatchCount++ }
$cell.Value = 1212
$cell = $ws.Cells.Item(23, 9)
if ($cell.Value2 -ne 1031) { $mismatchCount++ }
$cell.Value = 1212
$cell = $ws.Cells.Item(23, 10)
if ($cell.Value2 -ne 5999.5) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(23, 11)
if ($cell.Value2 -ne 1031) { $mismatchCount++ }
$cell.Value = 1212
$cell = $ws.Cells.Item(23, 12)
if ($cell.Value2 -ne 5999.5) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(23, 13)
if ($cell.Value2 -ne -748) { $mismatchCount++ }
$cell.Value = -929
$cell = $ws.Cells.Item(23, 14)
if ($cell.Value2 -ne -6565.5) { $mismatchCount++ }
$cell.Value = ""
# Row 29
$cell = $ws.Cells.Item(29, 8)
if ($cell.Value2 -ne 5500) { $mismatchCount++ }
$cell.Value = 505
$cell = $ws.Cells.Item(29, 9)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 505
$cell = $ws.Cells.Item(29, 10)
if ($cell.Value2 -ne 5500) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(29, 11)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 505
$cell = $ws.Cells.Item(29, 12)
if ($cell.Value2 -ne 5500) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(29, 13)
if ($cell.Value2 -ne $null -and $cell.Value2 -ne "") { $mismatchCount++ }
$cell.Value = -216
$cell = $ws.Cells.Item(29, 14)
if ($cell.Value2 -ne -6078) { $mismatchCount++ }
$cell.Value = ""
# Row 31
$cell = $ws.Cells.Item(31, 8)
if ($cell.Value2 -ne 2999) { $mismatchCount++ }
$cell.Value = 15000
$cell = $ws.Cells.Item(31, 9)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 15000
$cell = $ws.Cells.Item(31, 10)
if ($cell.Value2 -ne 2999) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(31, 11)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 15000
$cell = $ws.Cells.Item(31, 12)
if ($cell.Value2 -ne 2999) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(31, 13)
if ($cell.Value2 -ne $null -and $cell.Value2 -ne "") { $mismatchCount++ }
$cell.Value = -14748
$cell = $ws.Cells.Item(31, 14)
if ($cell.Value2 -ne -3503) { $mismatchCount++ }
$cell.Value = ""
# Row 44
$cell = $ws.Cells.Item(44, 8)
if ($cell.Value2 -ne 48413.145) { $mismatchCount++ }
$cell.Value = 38397.43
$cell = $ws.Cells.Item(44, 10)
if ($cell.Value2 -ne 48413.145) { $mismatchCount++ }
$cell.Value = 38397.43
$cell = $ws.Cells.Item(44, 12)
if ($cell.Value2 -ne 48413.145) { $mismatchCount++ }
$cell.Value = 38397.43
$cell = $ws.Cells.Item(44, 14)
if ($cell.Value2 -ne -49407.145) { $mismatchCount++ }
$cell.Value = -39391.43
# Row 99
$cell = $ws.Cells.Item(99, 8)
if ($cell.Value2 -ne 2778.611) { $mismatchCount++ }
$cell.Value = 2463.8235
$cell = $ws.Cells.Item(99, 9)
if ($cell.Value2 -ne 2951.125) { $mismatchCount++ }
$cell.Value = 2707.7144
$cell = $ws.Cells.Item(99, 10)
if ($cell.Value2 -ne 1398.5) { $mismatchCount++ }
$cell.Value = 1325.6666
$cell = $ws.Cells.Item(99, 11)
if ($cell.Value2 -ne 2951.125) { $mismatchCount++ }
$cell.Value = 2707.7144
$cell = $ws.Cells.Item(99, 12)
if ($cell.Value2 -ne 1398.5) { $mismatchCount++ }
$cell.Value = 1325.6666
$cell = $ws.Cells.Item(99, 13)
if ($cell.Value2 -ne -1453.125) { $mismatchCount++ }
$cell.Value = -1209.7144
$cell = $ws.Cells.Item(99, 14)
if ($cell.Value2 -ne -4394.5) { $mismatchCount++ }
$cell.Value = -4321.6666
# Row 107
$cell = $ws.Cells.Item(107, 8)
if ($cell.Value2 -ne 2411.889) { $mismatchCount++ }
$cell.Value = 2463.5
$cell = $ws.Cells.Item(107, 9)
if ($cell.Value2 -ne 2118) { $mismatchCount++ }
$cell.Value = 2141.8
$cell = $ws.Cells.Item(107, 11)
if ($cell.Value2 -ne 2118) { $mismatchCount++ }
$cell.Value = 2141.8
$cell = $ws.Cells.Item(107, 13)
if ($cell.Value2 -ne -198) { $mismatchCount++ }
$cell.Value = -221.8000000000002
# Row 134
$cell = $ws.Cells.Item(134, 8)
if ($cell.Value2 -ne 8749.454) { $mismatchCount++ }
$cell.Value = 8514.735000000001
$cell = $ws.Cells.Item(134, 9)
if ($cell.Value2 -ne 5420.0415) { $mismatchCount++ }
$cell.Value = 5234
$cell = $ws.Cells.Item(134, 11)
if ($cell.Value2 -ne 16260.1245) { $mismatchCount++ }
$cell.Value = 15702
$cell = $ws.Cells.Item(134, 13)
if ($cell.Value2 -ne -13725.1245) { $mismatchCount++ }
$cell.Value = -13167

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 5
$cell = $ws.Cells.Item(5, 8)
if ($cell.Value2 -ne 1307.25) { $mismatchCount++ }
$cell.Value = 756.55554
$cell = $ws.Cells.Item(5, 9)
if ($cell.Value2 -ne 885.6667) { $mismatchCount++ }
$cell.Value = 151.5
$cell = $ws.Cells.Item(5, 10)
if ($cell.Value2 -ne 1560.2) { $mismatchCount++ }
$cell.Value = 1966.6666
$cell = $ws.Cells.Item(5, 11)
if ($cell.Value2 -ne 885.6667) { $mismatchCount++ }
$cell.Value = 151.5
$cell = $ws.Cells.Item(5, 12)
if ($cell.Value2 -ne 1560.2) { $mismatchCount++ }
$cell.Value = 1966.6666
$cell = $ws.Cells.Item(5, 13)
if ($cell.Value2 -ne -773.6667) { $mismatchCount++ }
$cell.Value = -39.5
$cell = $ws.Cells.Item(5, 14)
if ($cell.Value2 -ne -1784.2) { $mismatchCount++ }
$cell.Value = -2190.6666
# Row 6
$cell = $ws.Cells.Item(6, 8)
if ($cell.Value2 -ne 12666.583) { $mismatchCount++ }
$cell.Value = 9668
$cell = $ws.Cells.Item(6, 9)
if ($cell.Value2 -ne 16475) { $mismatchCount++ }
$cell.Value = 12448.4
$cell = $ws.Cells.Item(6, 10)
if ($cell.Value2 -ne 5049.75) { $mismatchCount++ }
$cell.Value = 400
$cell = $ws.Cells.Item(6, 11)
if ($cell.Value2 -ne 16475) { $mismatchCount++ }
$cell.Value = 12448.4
$cell = $ws.Cells.Item(6, 12)
if ($cell.Value2 -ne 5049.75) { $mismatchCount++ }
$cell.Value = 400
$cell = $ws.Cells.Item(6, 13)
if ($cell.Value2 -ne -16362) { $mismatchCount++ }
$cell.Value = -12335.4
$cell = $ws.Cells.Item(6, 14)
if ($cell.Value2 -ne -5275.75) { $mismatchCount++ }
$cell.Value = -626
# Row 8
$cell = $ws.Cells.Item(8, 8)
if ($cell.Value2 -ne 4379.2) { $mismatchCount++ }
$cell.Value = 3367.125
$cell = $ws.Cells.Item(8, 9)
if ($cell.Value2 -ne 650) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(8, 10)
if ($cell.Value2 -ne 5311.5) { $mismatchCount++ }
$cell.Value = 3367.125
$cell = $ws.Cells.Item(8, 11)
if ($cell.Value2 -ne 650) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(8, 12)
if ($cell.Value2 -ne 5311.5) { $mismatchCount++ }
$cell.Value = 3367.125
$cell = $ws.Cells.Item(8, 13)
if ($cell.Value2 -ne -510) { $mismatchCount++ }
$cell.Value = ""
$cell = $ws.Cells.Item(8, 14)
if ($cell.Value2 -ne -5591.5) { $mismatchCount++ }
$cell.Value = -3647.125
# Row 10
$cell = $ws.Cells.Item(10, 8)
if ($cell.Value2 -ne 2628.2856) { $mismatchCount++ }
$cell.Value = 3239.6
$cell = $ws.Cells.Item(10, 9)
if ($cell.Value2 -ne 833.3333) { $mismatchCount++ }
$cell.Value = 1999.5
$cell = $ws.Cells.Item(10, 10)
if ($cell.Value2 -ne 3974.5) { $mismatchCount++ }
$cell.Value = 4066.3333
$cell = $ws.Cells.Item(10, 11)
if ($cell.Value2 -ne 833.3333) { $mismatchCount++ }
$cell.Value = 1999.5
$cell = $ws.Cells.Item(10, 12)
if ($cell.Value2 -ne 3974.5) { $mismatchCount++ }
$cell.Value = 4066.3333
$cell = $ws.Cells.Item(10, 13)
if ($cell.Value2 -ne -694.3333) { $mismatchCount++ }
$cell.Value = -1860.5
$cell = $ws.Cells.Item(10, 14)
if ($cell.Value2 -ne -4252.5) { $mismatchCount++ }
$cell.Value = -4344.3333
# Row 135
$cell = $ws.Cells.Item(135, 8)
if ($cell.Value2 -ne 99998.3) { $mismatchCount++ }
$cell.Value = 99999
$cell = $ws.Cells.Item(135, 10)
if ($cell.Value2 -ne 99998.25) { $mismatchCount++ }
$cell.Value = 99999
$cell = $ws.Cells.Item(135, 12)
if ($cell.Value2 -ne 99998.25) { $mismatchCount++ }
$cell.Value = 99999
$cell = $ws.Cells.Item(135, 14)
if ($cell.Value2 -ne -110138.25) { $mismatchCount++ }
$cell.Value = -110139

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 29
$cell = $ws.Cells.Item(29, 8)
if ($cell.Value2 -ne 322.68182) { $mismatchCount++ }
$cell.Value = 321.81818
$cell = $ws.Cells.Item(29, 10)
if ($cell.Value2 -ne 322.22223) { $mismatchCount++ }
$cell.Value = 321.16666
$cell = $ws.Cells.Item(29, 12)
if ($cell.Value2 -ne 966.66669) { $mismatchCount++ }
$cell.Value = 963.4999799999999
$cell = $ws.Cells.Item(29, 14)
if ($cell.Value2 -ne -1520.66669) { $mismatchCount++ }
$cell.Value = -1517.49998
# Row 43
$cell = $ws.Cells.Item(43, 8)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 20
$cell = $ws.Cells.Item(43, 9)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 20
$cell = $ws.Cells.Item(43, 11)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 60
$cell = $ws.Cells.Item(43, 13)
if ($cell.Value2 -ne $null -and $cell.Value2 -ne "") { $mismatchCount++ }
$cell.Value = 54

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$cell = $ws.Cells.Item(9, 8)
if ($cell.Value2 -ne 3526.6667) { $mismatchCount++ }
$cell.Value = 1317.8334
$cell = $ws.Cells.Item(9, 9)
if ($cell.Value2 -ne 230.8) { $mismatchCount++ }
$cell.Value = 226.75
$cell = $ws.Cells.Item(9, 10)
if ($cell.Value2 -ne 20006) { $mismatchCount++ }
$cell.Value = 3500
$cell = $ws.Cells.Item(9, 11)
if ($cell.Value2 -ne 230.8) { $mismatchCount++ }
$cell.Value = 226.75
$cell = $ws.Cells.Item(9, 12)
if ($cell.Value2 -ne 20006) { $mismatchCount++ }
$cell.Value = 3500
$cell = $ws.Cells.Item(9, 13)
if ($cell.Value2 -ne -60.80000000000001) { $mismatchCount++ }
$cell.Value = -56.75
$cell = $ws.Cells.Item(9, 14)
if ($cell.Value2 -ne -20346) { $mismatchCount++ }
$cell.Value = -3840
# Row 11
$cell = $ws.Cells.Item(11, 8)
if ($cell.Value2 -ne 1211115.2) { $mismatchCount++ }
$cell.Value = 1062566.5
$cell = $ws.Cells.Item(11, 9)
if ($cell.Value2 -ne 401542.72) { $mismatchCount++ }
$cell.Value = 480685.56
$cell = $ws.Cells.Item(11, 10)
if ($cell.Value2 -ne 2155616.5) { $mismatchCount++ }
$cell.Value = 1741427.5
$cell = $ws.Cells.Item(11, 11)
if ($cell.Value2 -ne 401542.72) { $mismatchCount++ }
$cell.Value = 480685.56
$cell = $ws.Cells.Item(11, 12)
if ($cell.Value2 -ne 2155616.5) { $mismatchCount++ }
$cell.Value = 1741427.5
$cell = $ws.Cells.Item(11, 13)
if ($cell.Value2 -ne -401403.72) { $mismatchCount++ }
$cell.Value = -480546.56
$cell = $ws.Cells.Item(11, 14)
if ($cell.Value2 -ne -2155894.5) { $mismatchCount++ }
$cell.Value = -1741705.5
# Row 12
$cell = $ws.Cells.Item(12, 8)
if ($cell.Value2 -ne 16175.5) { $mismatchCount++ }
$cell.Value = 15414.143
$cell = $ws.Cells.Item(12, 10)
if ($cell.Value2 -ne 19876) { $mismatchCount++ }
$cell.Value = 19333
$cell = $ws.Cells.Item(12, 12)
if ($cell.Value2 -ne 19876) { $mismatchCount++ }
$cell.Value = 19333
$cell = $ws.Cells.Item(12, 14)
if ($cell.Value2 -ne -20156) { $mismatchCount++ }
$cell.Value = -19613
# Row 25
$cell = $ws.Cells.Item(25, 8)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 9502.333000000001
$cell = $ws.Cells.Item(25, 9)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 3500
$cell = $ws.Cells.Item(25, 10)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 12503.5
$cell = $ws.Cells.Item(25, 11)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 3500
$cell = $ws.Cells.Item(25, 12)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 12503.5
$cell = $ws.Cells.Item(25, 13)
if ($cell.Value2 -ne $null -and $cell.Value2 -ne "") { $mismatchCount++ }
$cell.Value = -2971
$cell = $ws.Cells.Item(25, 14)
if ($cell.Value2 -ne $null -and $cell.Value2 -ne "") { $mismatchCount++ }
$cell.Value = -13561.5
# Row 29
$cell = $ws.Cells.Item(29, 8)
if ($cell.Value2 -ne 9539.799999999999) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(29, 10)
if ($cell.Value2 -ne 9539.799999999999) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(29, 12)
if ($cell.Value2 -ne 9539.799999999999) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(29, 14)
if ($cell.Value2 -ne -10119.8) { $mismatchCount++ }
$cell.Value = ""
# Row 132
$cell = $ws.Cells.Item(132, 8)
if ($cell.Value2 -ne 9865.1) { $mismatchCount++ }
$cell.Value = 9585.517
$cell = $ws.Cells.Item(132, 10)
if ($cell.Value2 -ne 6202.6) { $mismatchCount++ }
$cell.Value = 5368.5
$cell = $ws.Cells.Item(132, 12)
if ($cell.Value2 -ne 18607.8) { $mismatchCount++ }
$cell.Value = 16105.5
$cell = $ws.Cells.Item(132, 14)
if ($cell.Value2 -ne -23667.8) { $mismatchCount++ }
$cell.Value = -21165.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 44
$cell = $ws.Cells.Item(44, 8)
if ($cell.Value2 -ne 4250) { $mismatchCount++ }
$cell.Value = 4000
$cell = $ws.Cells.Item(44, 10)
if ($cell.Value2 -ne 4250) { $mismatchCount++ }
$cell.Value = 4000
$cell = $ws.Cells.Item(44, 12)
if ($cell.Value2 -ne 4250) { $mismatchCount++ }
$cell.Value = 4000
$cell = $ws.Cells.Item(44, 14)
if ($cell.Value2 -ne -5162) { $mismatchCount++ }
$cell.Value = -4912
# Row 50
$cell = $ws.Cells.Item(50, 8)
if ($cell.Value2 -ne 26146.666) { $mismatchCount++ }
$cell.Value = 27976.2
$cell = $ws.Cells.Item(50, 10)
if ($cell.Value2 -ne 27379.4) { $mismatchCount++ }
$cell.Value = 29974.5
$cell = $ws.Cells.Item(50, 12)
if ($cell.Value2 -ne 27379.4) { $mismatchCount++ }
$cell.Value = 29974.5
$cell = $ws.Cells.Item(50, 14)
if ($cell.Value2 -ne -28653.4) { $mismatchCount++ }
$cell.Value = -31248.5
# Row 57
$cell = $ws.Cells.Item(57, 8)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 4996
$cell = $ws.Cells.Item(57, 9)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 4996
$cell = $ws.Cells.Item(57, 11)
if ($cell.Value2 -ne 0) { $mismatchCount++ }
$cell.Value = 4996
$cell = $ws.Cells.Item(57, 13)
if ($cell.Value2 -ne $null -and $cell.Value2 -ne "") { $mismatchCount++ }
$cell.Value = -4430
# Row 132
$cell = $ws.Cells.Item(132, 8)
if ($cell.Value2 -ne 5537.25) { $mismatchCount++ }
$cell.Value = 5316.431
$cell = $ws.Cells.Item(132, 9)
if ($cell.Value2 -ne 5367.175) { $mismatchCount++ }
$cell.Value = 5200.881
$cell = $ws.Cells.Item(132, 10)
if ($cell.Value2 -ne 6387.625) { $mismatchCount++ }
$cell.Value = 5855.6665
$cell = $ws.Cells.Item(132, 11)
if ($cell.Value2 -ne 16101.525) { $mismatchCount++ }
$cell.Value = 15602.643
$cell = $ws.Cells.Item(132, 12)
if ($cell.Value2 -ne 19162.875) { $mismatchCount++ }
$cell.Value = 17566.9995
$cell = $ws.Cells.Item(132, 13)
if ($cell.Value2 -ne -13571.525) { $mismatchCount++ }
$cell.Value = -13072.643
$cell = $ws.Cells.Item(132, 14)
if ($cell.Value2 -ne -24222.875) { $mismatchCount++ }
$cell.Value = -22626.9995

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 4
$cell = $ws.Cells.Item(4, 8)
if ($cell.Value2 -ne 889.8) { $mismatchCount++ }
$cell.Value = 840.7143
$cell = $ws.Cells.Item(4, 9)
if ($cell.Value2 -ne 900) { $mismatchCount++ }
$cell.Value = 821.25
$cell = $ws.Cells.Item(4, 10)
if ($cell.Value2 -ne 887.25) { $mismatchCount++ }
$cell.Value = 866.6667
$cell = $ws.Cells.Item(4, 11)
if ($cell.Value2 -ne 900) { $mismatchCount++ }
$cell.Value = 821.25
$cell = $ws.Cells.Item(4, 12)
if ($cell.Value2 -ne 887.25) { $mismatchCount++ }
$cell.Value = 866.6667
$cell = $ws.Cells.Item(4, 13)
if ($cell.Value2 -ne -787) { $mismatchCount++ }
$cell.Value = -708.25
$cell = $ws.Cells.Item(4, 14)
if ($cell.Value2 -ne -1113.25) { $mismatchCount++ }
$cell.Value = -1092.6667
# Row 7
$cell = $ws.Cells.Item(7, 8)
if ($cell.Value2 -ne 585.6667) { $mismatchCount++ }
$cell.Value = 433.33334
$cell = $ws.Cells.Item(7, 9)
if ($cell.Value2 -ne 205) { $mismatchCount++ }
$cell.Value = 200
$cell = $ws.Cells.Item(7, 10)
if ($cell.Value2 -ne 966.3333) { $mismatchCount++ }
$cell.Value = 900
$cell = $ws.Cells.Item(7, 11)
if ($cell.Value2 -ne 205) { $mismatchCount++ }
$cell.Value = 200
$cell = $ws.Cells.Item(7, 12)
if ($cell.Value2 -ne 966.3333) { $mismatchCount++ }
$cell.Value = 900
$cell = $ws.Cells.Item(7, 13)
if ($cell.Value2 -ne -92) { $mismatchCount++ }
$cell.Value = -87
$cell = $ws.Cells.Item(7, 14)
if ($cell.Value2 -ne -1192.3333) { $mismatchCount++ }
$cell.Value = -1126
# Row 13
$cell = $ws.Cells.Item(13, 8)
if ($cell.Value2 -ne 4912) { $mismatchCount++ }
$cell.Value = 10973.5
$cell = $ws.Cells.Item(13, 9)
if ($cell.Value2 -ne 4861.3335) { $mismatchCount++ }
$cell.Value = 6005
$cell = $ws.Cells.Item(13, 10)
if ($cell.Value2 -ne 4962.6665) { $mismatchCount++ }
$cell.Value = 12629.667
$cell = $ws.Cells.Item(13, 11)
if ($cell.Value2 -ne 4861.3335) { $mismatchCount++ }
$cell.Value = 6005
$cell = $ws.Cells.Item(13, 12)
if ($cell.Value2 -ne 4962.6665) { $mismatchCount++ }
$cell.Value = 12629.667
$cell = $ws.Cells.Item(13, 13)
if ($cell.Value2 -ne -4721.3335) { $mismatchCount++ }
$cell.Value = -5865
$cell = $ws.Cells.Item(13, 14)
if ($cell.Value2 -ne -5242.6665) { $mismatchCount++ }
$cell.Value = -12909.667
# Row 17
$cell = $ws.Cells.Item(17, 8)
if ($cell.Value2 -ne 6243.5713) { $mismatchCount++ }
$cell.Value = 11114507
$cell = $ws.Cells.Item(17, 9)
if ($cell.Value2 -ne 4781.4) { $mismatchCount++ }
$cell.Value = 3820.875
$cell = $ws.Cells.Item(17, 10)
if ($cell.Value2 -ne 9899) { $mismatchCount++ }
$cell.Value = 100000000
$cell = $ws.Cells.Item(17, 11)
if ($cell.Value2 -ne 4781.4) { $mismatchCount++ }
$cell.Value = 3820.875
$cell = $ws.Cells.Item(17, 12)
if ($cell.Value2 -ne 9899) { $mismatchCount++ }
$cell.Value = 100000000
$cell = $ws.Cells.Item(17, 13)
if ($cell.Value2 -ne -4609.4) { $mismatchCount++ }
$cell.Value = -3648.875
$cell = $ws.Cells.Item(17, 14)
if ($cell.Value2 -ne -10243) { $mismatchCount++ }
$cell.Value = -100000344
# Row 23
$cell = $ws.Cells.Item(23, 8)
if ($cell.Value2 -ne 3530) { $mismatchCount++ }
$cell.Value = 512.875
$cell = $ws.Cells.Item(23, 9)
if ($cell.Value2 -ne 285) { $mismatchCount++ }
$cell.Value = 250.33333
$cell = $ws.Cells.Item(23, 10)
if ($cell.Value2 -ne 6775) { $mismatchCount++ }
$cell.Value = 1300.5
$cell = $ws.Cells.Item(23, 11)
if ($cell.Value2 -ne 285) { $mismatchCount++ }
$cell.Value = 250.33333
$cell = $ws.Cells.Item(23, 12)
if ($cell.Value2 -ne 6775) { $mismatchCount++ }
$cell.Value = 1300.5
$cell = $ws.Cells.Item(23, 13)
if ($cell.Value2 -ne -56) { $mismatchCount++ }
$cell.Value = -21.33332999999999
$cell = $ws.Cells.Item(23, 14)
if ($cell.Value2 -ne -7233) { $mismatchCount++ }
$cell.Value = -1758.5
# Row 34
$cell = $ws.Cells.Item(34, 8)
if ($cell.Value2 -ne 20007.666) { $mismatchCount++ }
$cell.Value = 20026
$cell = $ws.Cells.Item(34, 9)
if ($cell.Value2 -ne 20012.5) { $mismatchCount++ }
$cell.Value = 20026
$cell = $ws.Cells.Item(34, 10)
if ($cell.Value2 -ne 19998) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(34, 11)
if ($cell.Value2 -ne 20012.5) { $mismatchCount++ }
$cell.Value = 20026
$cell = $ws.Cells.Item(34, 12)
if ($cell.Value2 -ne 19998) { $mismatchCount++ }
$cell.Value = 0
$cell = $ws.Cells.Item(34, 13)
if ($cell.Value2 -ne -19809.5) { $mismatchCount++ }
$cell.Value = -19823
$cell = $ws.Cells.Item(34, 14)
if ($cell.Value2 -ne -20404) { $mismatchCount++ }
$cell.Value = ""
# Row 37
$cell = $ws.Cells.Item(37, 8)
if ($cell.Value2 -ne 66664.336) { $mismatchCount++ }
$cell.Value = 61997.25
$cell = $ws.Cells.Item(37, 10)
if ($cell.Value2 -ne 99994) { $mismatchCount++ }
$cell.Value = 73995
$cell = $ws.Cells.Item(37, 12)
if ($cell.Value2 -ne 99994) { $mismatchCount++ }
$cell.Value = 73995
$cell = $ws.Cells.Item(37, 14)
if ($cell.Value2 -ne -100400) { $mismatchCount++ }
$cell.Value = -74401
# Row 50
$cell = $ws.Cells.Item(50, 8)
if ($cell.Value2 -ne 30440) { $mismatchCount++ }
$cell.Value = 27882
$cell = $ws.Cells.Item(50, 10)
if ($cell.Value2 -ne 30440) { $mismatchCount++ }
$cell.Value = 27882
$cell = $ws.Cells.Item(50, 12)
if ($cell.Value2 -ne 30440) { $mismatchCount++ }
$cell.Value = 27882
$cell = $ws.Cells.Item(50, 14)
if ($cell.Value2 -ne -31702) { $mismatchCount++ }
$cell.Value = -29144
# Row 55
$cell = $ws.Cells.Item(55, 8)
if ($cell.Value2 -ne 6283.4287) { $mismatchCount++ }
$cell.Value = 4048.25
$cell = $ws.Cells.Item(55, 9)
if ($cell.Value2 -ne 2799.2) { $mismatchCount++ }
$cell.Value = 2485.1428
$cell = $ws.Cells.Item(55, 10)
if ($cell.Value2 -ne 14994) { $mismatchCount++ }
$cell.Value = 14990
$cell = $ws.Cells.Item(55, 11)
if ($cell.Value2 -ne 2799.2) { $mismatchCount++ }
$cell.Value = 2485.1428
$cell = $ws.Cells.Item(55, 12)
if ($cell.Value2 -ne 14994) { $mismatchCount++ }
$cell.Value = 14990
$cell = $ws.Cells.Item(55, 13)
if ($cell.Value2 -ne -2522.2) { $mismatchCount++ }
$cell.Value = -2208.1428
$cell = $ws.Cells.Item(55, 14)
if ($cell.Value2 -ne -15548) { $mismatchCount++ }
$cell.Value = -15544

if ($mismatchCount -gt 0) {
    Write-Output ("Warning: " + $mismatchCount + " cell(s) did not match the expected prior value before being overwritten.")
} else {
    Write-Output "All prior values matched expectations; update applied."
}
